$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 448.5567906052825
$ws.Range("C2").Value = 2282.404201474098
$ws.Range("D2").Value = 3956.640472599907
